$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 22:16"

# Row 6: Estados Unidos
$ws.Range("A6").Value = "Estados Unidos"
$ws.Range("B6").Value = 64769
$ws.Range("C6").Value = 9913
$ws.Range("D6").Value = 393
$ws.Range("E6").Value = 63466
$ws.Range("F6").Value = 1411
$ws.Range("G6").Value = 130
$ws.Range("H6").Value = 910

# Row 17: Canada
$ws.Range("A17").Value = "Canada"
$ws.Range("B17").Value = 3306
$ws.Range("C17").Value = 514
$ws.Range("D17").Value = 185
$ws.Range("E17").Value = 3091
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 30

# Row 22: Turquia
$ws.Range("A22").Value = "Turquia"
$ws.Range("B22").Value = 2433
$ws.Range("C22").Value = 561
$ws.Range("D22").Value = 26
$ws.Range("E22").Value = 2348
$ws.Range("F22").Value = 136
$ws.Range("G22").Value = 15
$ws.Range("H22").Value = 59

# Row 23: Australia
$ws.Range("A23").Value = "Australia"
$ws.Range("B23").Value = 2431
$ws.Range("C23").Value = 114
$ws.Range("D23").Value = 118
$ws.Range("E23").Value = 2304
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 9

# Row 24: Israel
$ws.Range("A24").Value = "Israel"
$ws.Range("B24").Value = 2369
$ws.Range("C24").Value = 439
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = 2306
$ws.Range("F24").Value = 37
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 5

# Row 74: San Marino
$ws.Range("A74").Value = "San Marino"
$ws.Range("B74").Value = 208
$ws.Range("C74").Value = 21
$ws.Range("D74").Value = 4
$ws.Range("E74").Value = 183
$ws.Range("F74").Value = 12
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 21

# Row 75: Nueva Zelanda
$ws.Range("A75").Value = "Nueva Zelanda"
$ws.Range("B75").Value = 205
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 22
$ws.Range("E75").Value = 183
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0

# Row 76: Costa Rica
$ws.Range("A76").Value = "Costa Rica"
$ws.Range("B76").Value = 201
$ws.Range("C76").Value = 24
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 197
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 2

# Row 77: Kuwait
$ws.Range("A77").Value = "Kuwait"
$ws.Range("B77").Value = 195
$ws.Range("C77").Value = 4
$ws.Range("D77").Value = 43
$ws.Range("E77").Value = 152
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0

# Row 78: Uruguay
$ws.Range("A78").Value = "Uruguay"
$ws.Range("B78").Value = 189
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 189
$ws.Range("F78").Value = 3
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 0

# Row 79: Principado de Andorra
$ws.Range("A79").Value = "Principado de Andorra"
$ws.Range("B79").Value = 188
$ws.Range("C79").Value = 24
$ws.Range("D79").Value = 1
$ws.Range("E79").Value = 186
$ws.Range("F79").Value = 6
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 1

# Row 87: Ucrania
$ws.Range("A87").Value = "Ucrania"
$ws.Range("B87").Value = 145
$ws.Range("C87").Value = 43
$ws.Range("D87").Value = 1
$ws.Range("E87").Value = 139
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 5

# Row 88: Vietnam
$ws.Range("A88").Value = "Vietnam"
$ws.Range("B88").Value = 141
$ws.Range("C88").Value = 7
$ws.Range("D88").Value = 17
$ws.Range("E88").Value = 124
$ws.Range("F88").Value = 3
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0

# Row 89: Republica de Chipre
$ws.Range("A89").Value = "Republica de Chipre"
$ws.Range("B89").Value = 132
$ws.Range("C89").Value = 8
$ws.Range("D89").Value = 3
$ws.Range("E89").Value = 126
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 3

# Row 90: Islas Feroe
$ws.Range("A90").Value = "Islas Feroe"
$ws.Range("B90").Value = 132
$ws.Range("C90").Value = 10
$ws.Range("D90").Value = 38
$ws.Range("E90").Value = 94
$ws.Range("F90").Value = 2
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 0

# Row 91: Malta
$ws.Range("A91").Value = "Malta"
$ws.Range("B91").Value = 129
$ws.Range("C91").Value = 19
$ws.Range("D91").Value = 2
$ws.Range("E91").Value = 127
$ws.Range("F91").Value = 1
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0
